$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column-A labels (phi_ values), rows 2..22
$labels = @(
    "phi_-0.0000",
    "phi_0.0500",
    "phi_0.1000",
    "phi_0.1500",
    "phi_0.2000",
    "phi_0.2500",
    "phi_0.3000",
    "phi_0.3500",
    "phi_0.4000",
    "phi_0.4500",
    "phi_0.5000",
    "phi_0.5500",
    "phi_0.6000",
    "phi_0.6500",
    "phi_0.7000",
    "phi_0.7500",
    "phi_0.8000",
    "phi_0.8500",
    "phi_0.9000",
    "phi_0.9500",
    "phi_1.0000"
)

# New data for columns B (coarse), C (medium), D (fine), E (x)
$data = @(
    @(0.446113646,        0.442592859,        0.4470689,          -0.000000713393092),
    @(0.4246112354279316, 0.4262847582645582, 0.435861379980351,  0.04999932227656259),
    @(0.362232850563531,  0.3635276939148809, 0.3698628781433433, 0.09999935794621718),
    @(0.2869568017485726, 0.2929609431107805, 0.2979116755464833, 0.1499993936158718),
    @(0.2085413773140853, 0.2167929287941868, 0.2201327974124013, 0.1999994292855264),
    @(0.1351285381754221, 0.1423025141596007, 0.14372573426556,   0.249999464955181),
    @(0.07138259991756449,0.07674601964131285,0.0771815567248072, 0.2999995006248355),
    @(0.02297510738591806,0.02629345425254802,0.02552679883793104,0.3499995362944902),
    @(0.04411748410674251,0.04416111023660232,0.03993078715810641,0.3999995719641448),
    @(0.08988136754102584,0.08991195453941069,0.08629661294338958,0.4499996076337994),
    @(0.1361954070995965, 0.1359988540096508, 0.1326337387927119, 0.499999643303454),
    @(0.1840127752610901, 0.1821744756461599, 0.178864023292082,  0.5499996789731085),
    @(0.2383783516901995, 0.2323943639211221, 0.2293526279103272, 0.599999714642763),
    @(0.3127474316144384, 0.2995419889334269, 0.2986893882845401, 0.6499997503124176),
    @(0.4283596859988836, 0.4116351294545521, 0.4192823897277798, 0.6999997859820722),
    @(0.6075918816272678, 0.5950562549488467, 0.6147961707990083, 0.7499998216517269),
    @(0.8359285424689332, 0.8321582760263507, 0.8591015965897354, 0.7999998573213815),
    @(0.9494482100330378, 0.9767981609654574, 0.9882839937274183, 0.8499998929910361),
    @(0.9028568367635105, 0.9367215821998738, 0.944724354652655,  0.8999999286606907),
    @(0.7217200486092825, 0.7447144090249544, 0.7355765974985735, 0.9499999643303453),
    @(0.527168632,        0.546208382,        0.515793443,        1)
)

for ($i = 0; $i -lt 21; $i++) {
    $row = $i + 2

    $ws.Cells.Item($row, 1).Value = $labels[$i]
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
    $ws.Cells.Item($row, 5).Value = $data[$i][3]
}
